$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 47 keeps its own record; only the Ost/Nord coordinates get rounded
# to whole metres and the Starttid/Sluttid ("00:00") cells are dropped.
$ws.Range("Q47").Value = 439827
$ws.Range("R47").Value = 6952233
$ws.Range("Z47").ClearContents()
$ws.Range("AB47").ClearContents()

# Rows 48-51 are re-shuffled: each target row receives the species /
# find data that used to live in a different row (48<-50, 49<-51,
# 50<-48, 51<-49), the Ost/Nord values are rounded to whole metres,
# and the Starttid/Sluttid cells are dropped everywhere.

# Row 48 <- old row 50 data
$ws.Range("A48").Value = 111974191
$ws.Range("B48").Value = 90652
$ws.Range("E48").Value = 3100
$ws.Range("F48").Value = "Talltaggsvamp"
$ws.Range("G48").Value = "Bankera fuligineoalba"
$ws.Range("H48").Value = "(Schmidt : Fr.) Pouzar"
$ws.Range("Q48").Value = 439978
$ws.Range("R48").Value = 6952214
$ws.Range("Z48").ClearContents()
$ws.Range("AB48").ClearContents()
$ws.Range("AI48").Value = "äldre renbetad ristallskog med lavfläckar på torr moränmark"

# Row 49 <- old row 51 data
$ws.Range("A49").Value = 111974186
$ws.Range("B49").Value = 90682
$ws.Range("E49").Value = 2059
$ws.Range("F49").Value = "Skrovlig taggsvamp"
$ws.Range("G49").Value = "Hydnellum scabrosum"
$ws.Range("H49").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("Q49").Value = 439860
$ws.Range("R49").Value = 6952250
$ws.Range("Z49").ClearContents()
$ws.Range("AB49").ClearContents()

# Row 50 <- old row 48 data
$ws.Range("A50").Value = 111974187
$ws.Range("B50").Value = 90710
$ws.Range("E50").Value = 5449
$ws.Range("F50").Value = "Svart taggsvamp"
$ws.Range("G50").Value = "Phellodon niger"
$ws.Range("H50").Value = "(Fr.:Fr.) P.Karst."
$ws.Range("Q50").Value = 439865
$ws.Range("R50").Value = 6952242
$ws.Range("Z50").ClearContents()
$ws.Range("AB50").ClearContents()
$ws.Range("AI50").Value = "äldre renbetad ristallskog med lavfläckar på torr moränmark, under tallåga"

# Row 51 <- old row 49 data
$ws.Range("A51").Value = 111974188
$ws.Range("B51").Value = 90652
$ws.Range("E51").Value = 3100
$ws.Range("F51").Value = "Talltaggsvamp"
$ws.Range("G51").Value = "Bankera fuligineoalba"
$ws.Range("H51").Value = "(Schmidt : Fr.) Pouzar"
$ws.Range("Q51").Value = 439870
$ws.Range("R51").Value = 6952225
$ws.Range("Z51").ClearContents()
$ws.Range("AB51").ClearContents()
